$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 corresponds to 10-January-2025.
# Switch the day's status from "At Work" to "Sick Leave".
$ws.Range("C20").Value = 0
$ws.Range("E20").Value = 1
$ws.Range("H20").Value = "Sick Leave"
